# Updated cryptos list on Wed Nov  1 18:57:57 UTC 2023 with GitHub Actions
#
# The sheet holds one crypto-currency per row (coin name, link, price,
# 1h volume%). This refresh pass overwrites the Price/Volume columns with
# the latest scrape, and a handful of rows also swapped rank (coin name +
# link moved to a neighboring row) because the source ranking reordered.
#
# Every value below is written as literal text (matching the original
# inlineStr cells) - several "price" strings (e.g. "1.00", "0.100",
# "226.20") are valid numeric literals, and Excel's COM layer silently
# re-types a plain .Value assignment like that into a Double (dropping
# the trailing zero / thousands-look formatting). To avoid that we
# briefly force the cell to Text ("@") before the write, then restore its
# original style so the saved cell keeps no stray number-format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new literal text, in the same row order as the diff
$updates = @(
    @{ Cell = "D2";  Text = "34.593.23" }
    @{ Cell = "E2";  Text = "  +0.46%  " }

    @{ Cell = "D3";  Text = "1.819.34" }
    @{ Cell = "E3";  Text = "  +0.67%  " }

    @{ Cell = "E4";  Text = "  +0.22%  " }

    @{ Cell = "D5";  Text = "226.20" }
    @{ Cell = "E5";  Text = "  +0.17%  " }

    @{ Cell = "D6";  Text = "0.607" }
    @{ Cell = "E6";  Text = "  +0.86%  " }

    @{ Cell = "D8";  Text = "44.83" }
    @{ Cell = "E8";  Text = "  +23.49%  " }

    @{ Cell = "D9";  Text = "0.298" }
    @{ Cell = "E9";  Text = "  +1.52%  " }

    @{ Cell = "D10"; Text = "0.0678" }
    @{ Cell = "E10"; Text = "  -0.42%  " }

    @{ Cell = "D11"; Text = "0.100" }
    @{ Cell = "E11"; Text = "  +3.84%  " }

    @{ Cell = "D12"; Text = "2.079.32" }
    @{ Cell = "E12"; Text = "  +0.65%  " }

    @{ Cell = "D13"; Text = "1.821.39" }
    @{ Cell = "E13"; Text = "  +0.20%  " }

    @{ Cell = "D14"; Text = "11.14" }

    # Polygon / Polkadot swapped rank positions
    @{ Cell = "B15"; Text = "Polkadot" }
    @{ Cell = "C15"; Text = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" }
    @{ Cell = "D15"; Text = "4.55" }
    @{ Cell = "E15"; Text = "  +2.33%  " }

    @{ Cell = "B16"; Text = "Polygon" }
    @{ Cell = "C16"; Text = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" }
    @{ Cell = "D16"; Text = "0.643" }
    @{ Cell = "E16"; Text = "  +1.93%  " }

    @{ Cell = "D17"; Text = "34.542.54" }
    @{ Cell = "E17"; Text = "  +0.41%  " }

    @{ Cell = "D18"; Text = "68.07" }
    @{ Cell = "E18"; Text = "  -0.81%  " }

    @{ Cell = "D19"; Text = "241.85" }
    @{ Cell = "E19"; Text = "  -0.41%  " }

    @{ Cell = "D20"; Text = "0.0₃0781" }
    @{ Cell = "E20"; Text = "  +0.88%  " }

    @{ Cell = "D21"; Text = "11.68" }
    @{ Cell = "E21"; Text = "  +3.82%  " }

    # Dai / Uniswap swapped rank positions
    @{ Cell = "B22"; Text = "Uniswap" }
    @{ Cell = "C22"; Text = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" }
    @{ Cell = "D22"; Text = "4.58" }
    @{ Cell = "E22"; Text = "  +11.56%  " }

    @{ Cell = "B23"; Text = "Dai" }
    @{ Cell = "C23"; Text = "https://coinranking.com/coin/MoTuySvg7+dai-dai" }
    @{ Cell = "D23"; Text = "1.00" }
    @{ Cell = "E23"; Text = "  +0.10%  " }

    @{ Cell = "E24"; Text = "  -1.95%  " }

    @{ Cell = "D25"; Text = "170.89" }
    @{ Cell = "E25"; Text = "  -0.14%  " }

    @{ Cell = "D26"; Text = "7.82" }
    @{ Cell = "E26"; Text = "  -0.62%  " }

    @{ Cell = "D27"; Text = "17.73" }
    @{ Cell = "E27"; Text = "  +2.13%  " }

    @{ Cell = "E28"; Text = "  +0.04%  " }

    @{ Cell = "E29"; Text = "  +0.13%  " }

    @{ Cell = "D30"; Text = "3.86" }
    @{ Cell = "E30"; Text = "  +1.48%  " }

    @{ Cell = "E31"; Text = "  +1.26%  " }

    @{ Cell = "D32"; Text = "3.95" }
    @{ Cell = "E32"; Text = "  +0.71%  " }

    @{ Cell = "D33"; Text = "0.0521" }
    @{ Cell = "E33"; Text = "  +0.90%  " }

    @{ Cell = "D34"; Text = "1.84" }
    @{ Cell = "E34"; Text = "  +2.24%  " }

    @{ Cell = "D35"; Text = "89.33" }

    @{ Cell = "D36"; Text = "0.660" }
    @{ Cell = "E36"; Text = "  +1.18%  " }

    @{ Cell = "D37"; Text = "15.32" }
    @{ Cell = "E37"; Text = "  +14.40%  " }

    @{ Cell = "D38"; Text = "1.320.06" }
    @{ Cell = "E38"; Text = "  -3.07%  " }

    @{ Cell = "E39"; Text = "  +0.20%  " }

    @{ Cell = "D40"; Text = "2.41" }
    @{ Cell = "E40"; Text = "  +1.67%  " }

    @{ Cell = "D41"; Text = "0.0191" }
    @{ Cell = "E41"; Text = "  +2.53%  " }

    # WEMIXToken / ARBITRUM swapped rank positions
    @{ Cell = "B42"; Text = "ARBITRUM" }
    @{ Cell = "C42"; Text = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" }
    @{ Cell = "D42"; Text = "0.959" }
    @{ Cell = "E42"; Text = "  +2.06%  " }

    @{ Cell = "B43"; Text = "WEMIXToken" }
    @{ Cell = "C43"; Text = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" }
    @{ Cell = "D43"; Text = "1.22" }
    @{ Cell = "E43"; Text = "  +4.89%  " }

    @{ Cell = "D44"; Text = "2.42" }
    @{ Cell = "E44"; Text = "  -0.05%  " }

    @{ Cell = "E45"; Text = "  +0.87%  " }

    @{ Cell = "E46"; Text = "  +4.36%  " }

    @{ Cell = "D47"; Text = "1.979.66" }
    @{ Cell = "E47"; Text = "  +0.67%  " }

    @{ Cell = "D48"; Text = "5.89" }
    @{ Cell = "E48"; Text = "  +1.43%  " }

    @{ Cell = "E49"; Text = "  +0.08%  " }

    @{ Cell = "D50"; Text = "101.43" }
    @{ Cell = "E50"; Text = "  -1.08%  " }

    # Cronos -> NEARProtocol (rank slot replaced outright)
    @{ Cell = "B51"; Text = "NEARProtocol" }
    @{ Cell = "C51"; Text = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" }
    @{ Cell = "D51"; Text = "1.56" }
    @{ Cell = "E51"; Text = "  +19.29%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Text

    # Would a plain .Value assignment get auto-coerced into a Double by
    # Excel (e.g. "1.00" -> 1, "0.100" -> 0.1, "226.20" -> 226.2)? If so,
    # force Text format for the write, then restore the cell's original
    # style so no stray number-format is left behind.
    $isPlainNumber = $text.Trim() -match '^[+-]?(\d+\.?\d*|\.\d+)$'

    if ($isPlainNumber) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = $origStyle
    } else {
        $cell.Value = $text
    }
}
